# Update cryptos list data to reflect latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.254.49'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.559.53'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.59'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.556.56'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.78%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '8.11'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.137'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.411'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.163.00'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('E14').Value = '  -3.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.22'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.564.02'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.348.49'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.33'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.61%  '
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.92'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '429.51'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.605'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.88'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.698.86'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000121'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.25'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.65%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +0.83%  '
$ws.Range('E32').Value = '  -4.81%  '
$ws.Range('E33').Value = '  -4.16%  '
$ws.Range('B34').Value = 'RenzoRestakedETH'
$ws.Range('C34').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.553.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.42'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.75'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '174.59'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0859'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.35'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.895'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E45').Value = '  -6.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.63'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.02'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -7.59%  '
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('E50').Value = '  -3.74%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.943'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.26%  '
